# Automatic update of files.
# Updates the "Förändrad" (Changed) date in column C for every data row
# from serial 45181 (2023-09-12) to serial 45182 (2023-09-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C ("Förändrad")
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
